# Change the annotation delimiter sample data from comma-separated to
# double-backslash-separated (comma is used as the argument delimiter).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Embedded,Id"  -> "Embedded\\Id"
$ws.Range("C9").Value = "Embedded\\Id"

# "Property,Id"  -> "Property\\Id"
$ws.Range("E23").Value = "Property\\Id"

# Description / help text cells: "カンマ区切り" -> "2重のバックスラッシュ区切り"
$ws.Range("C10").Value = "アノテーション付きバリューオブジェクトのサンプル。アノテーションは2重のバックスラッシュ区切りで複数記述できます。このクラスは単にサンプルです。実際の動作には利用されません。"
$ws.Range("G23").Value = "アノテーションのテスト。アノテーションは2重のバックスラッシュ区切りで複数記述できます。"
$ws.Range("G24").Value = "アノテーションのテスト。アノテーションは2重のバックスラッシュ区切りで複数記述できます。"

# Move the active selection from G11 to C11 (matches the saved sheet view state).
$ws.Range("C11").Select()

# Best-effort: reposition the workbook window (cosmetic window-chrome state;
# harmless no-op if the host doesn't persist it).
try {
    $wb.Windows.Item(1).Left = 660
    $wb.Windows.Item(1).Top = 500
} catch {
}
